# Insert a new "Plane" column before the existing X column (column C),
# shifting X, Y, Z, X_optimized, Y_optimized, Z_optimized, Array_Channel,
# Stagebox, Stagebox_Channel one column to the right (D..L).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C:C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "Plane"

# Mic-name prefix -> plane name lookup, derived from the existing
# Stagebox values for each row group (B->right, C->front, D->left, A->top).
$planeByPrefix = @{
    "A" = "top"
    "B" = "right"
    "C" = "front"
    "D" = "left"
}

# Find the last used row on the sheet and fill in the Plane column for
# every data row based on the Mic_Index (column B) label's leading letter.
$lastRow = $ws.Cells.Item($ws.Rows.Count, "B").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $micName = $ws.Cells.Item($r, 2).Value()
    if ($micName) {
        $prefix = $micName.Substring(0, 1)
        $plane = $planeByPrefix[$prefix]
        $ws.Cells.Item($r, 3).Value = $plane
    }
}
